$p = $ppt.ActivePresentation

# The deck currently ends with the "Thank You" slide (slide 5). Insert a
# new "GitHub Repository/Code" slide right before it, using the same
# "Title and Content" layout as the Thank You slide, so the deck becomes:
# 1..4 (unchanged), 5 = GitHub Repository/Code (new), 6 = Thank You.
$thankYou = $p.Slides.Item(5)
$layout = $thankYou.CustomLayout

$s = $p.Slides.AddSlide(5, $layout)

# Title placeholder
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "GitHub Repository/Code"
$title.LanguageID = "en-IN"

# Content placeholder: text + exact position/size (EMU 929640,2467293,10515600,683895)
$content = $s.Shapes.Item(2)
$contentRange = $content.TextFrame.TextRange
$contentRange.Text = "https://github.com/AdityaViswa/Stochastic-Path-Finding"
$contentRange.LanguageID = "en-IN"

$content.Left = 73.20000457763672
$content.Top = 194.27503967285156
$content.Width = 828.0
$content.Height = 53.85000228881836
